$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simpangan Baku Standard Deviasi")
$ws.Range("H7").Value = "S2  ="
$chars = $ws.Range("H7").Characters(2,2)
$chars.Font.Superscript = $true
